# Update column F (dSF) values for several rows as part of a data repull /
# push-all-data / mean-calculation refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -5
    4  = -2
    5  = -2
    6  = -2
    7  = 0
    8  = -5
    9  = -2
    10 = -2
    11 = -2
    14 = -2
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = -8
    21 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
